$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CushionDB header to include BEARINGS
$ws.Range("D1").Value = "CushionDB,BEARINGS"

# Add the new BEARINGS rule syntax to row 4 (Rule 'N3')
$ws.Range("D4").Value = "BEARINGS=['A','B']"

# Move the active selection from D3 to D4
$ws.Range("D4").Select()
